$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mobility")

# Add new row 4 data: variable name, value, and description
$ws.Range("A4").Value = "net_to_gross_factor"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Factor to convert (aproximately) net tons to gross tons, based on a full train operation (coef)."

# Update the active selection to B4
$ws.Range("B4").Select()
